# Rename the worksheet to unify DataNode / DataTable / Entity naming
# (was "Property1", matching the author's other Excel_Ini sheets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Restore the author's last on-disk selection/view state for the sheet.
$ws.Range("C41").Select() | Out-Null
